# Update sample fastq file paths to dev bucket gcloud location instead of SRA
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# --- Replace SRA accession values with gs:// bucket paths (fastq_file_1 / fastq_file_2) ---
$ws.Range("O4").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579104_1.fastq.gz"
$ws.Range("O5").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579105_1.fastq.gz"
$ws.Range("O6").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579106_1.fastq.gz"
$ws.Range("O7").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579107_1.fastq.gz"
$ws.Range("O8").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579108_1.fastq.gz"
$ws.Range("O9").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579109_1.fastq.gz"

$ws.Range("P4").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579104_2.fastq.gz"
$ws.Range("P5").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579105_2.fastq.gz"
$ws.Range("P6").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579106_2.fastq.gz"
$ws.Range("P7").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579107_2.fastq.gz"
$ws.Range("P8").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579108_2.fastq.gz"
$ws.Range("P9").Value = "gs://chips2-test-data/atac-test-sample-fastq/SRR16579109_2.fastq.gz"

# --- Widen column P (fastq_file_2) to fit the longer gs:// paths ---
$ws.Columns.Item(16).ColumnWidth = 56.333333333

# --- Make "metadata" the active sheet (was "workflow_config") and move the cursor ---
$ws.Activate() | Out-Null
$ws.Range("Q18").Select() | Out-Null
